# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (fund-level holdings) positioned
#    right before the "总计" (totals) sheet.
# 2) Insert a new top row into "总计" with the 2022-Q1 aggregate figures,
#    shifting the existing rows down and renumbering the index column.
#
# NOTE: worksheet object handles in this host are positional snapshots —
# they go stale after any operation that adds/moves/removes a sheet (or
# inserts/deletes rows). So after such an operation we always re-fetch
# the worksheet we need via $wb.Worksheets.Item(<name>) instead of
# reusing an old variable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add()
$new.Name = "2022-Q1"

$totalsForMove = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Item("2022-Q1")
$ws.Move($totalsForMove)

# Re-fetch (Move invalidates prior handles) and pull header / index-column
# formatting (style "2" in the source sheets) from an existing quarter
# sheet so the new sheet matches the look of its siblings. A single cell
# is copied and then "tiled" across the destination range by PasteSpecial.
$ws = $wb.Worksheets.Item("2022-Q1")
$template = $wb.Worksheets.Item("2021-Q4")

$template.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)

# Header row
$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# Columns B-G hold text (fund code / name / percentages kept as strings);
# a leading apostrophe forces text storage so values like "0.10" keep
# their trailing zero and "002685" keeps its leading zero. ClearFormats
# afterwards drops the "stored as text" quote-prefix marker Excel adds,
# so the cells end up with no explicit style (matching the sibling sheets).
$rows = @(
    @(0,"002685","中欧丰泓沪港深灵活配置混合A","59.37","94.70","4.48","2.6598",7),
    @(1,"010671","景顺长城大中华混合(QDII)美元","10.35","82.59","4.97","0.5144",7),
    @(2,"262001","景顺长城大中华混合(QDII)","10.35","82.59","4.97","0.5144",7),
    @(3,"002686","中欧丰泓沪港深灵活配置混合C","7.65","94.70","4.48","0.3427",7),
    @(4,"519779","交银施罗德沪港深价值精选灵活配置混合","5.13","84.44","5.26","0.2698",2),
    @(5,"006202","交银施罗德核心资产混合","0.72","79.85","5.45","0.0392",1),
    @(6,"004532","民生加银中证港股通高股息精选指数A","0.26","94.88","4.48","0.0116",5),
    @(7,"004533","民生加银中证港股通高股息精选指数C","0.10","94.88","4.48","0.0045",5),
    @(8,"011647","博时港股通红利精选混合A","0.13","92.10","3.28","0.0043",8),
    @(9,"011648","博时港股通红利精选混合C","0.02","92.10","3.28","0.0007",8)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = "'" + $row[1]
    $ws.Cells.Item($r,3).Value = "'" + $row[2]
    $ws.Cells.Item($r,4).Value = "'" + $row[3]
    $ws.Cells.Item($r,5).Value = "'" + $row[4]
    $ws.Cells.Item($r,6).Value = "'" + $row[5]
    $ws.Cells.Item($r,7).Value = "'" + $row[6]
    $ws.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

$ws.Range("B2:G11").ClearFormats()

# ---------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to "总计"
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()

# Re-fetch after the row insert and clean up the inherited blank formatting.
$totals = $wb.Worksheets.Item("总计")
$totals.Range("A2:D2").ClearFormats()

$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(2,2).Value = "2022-Q1"
$totals.Cells.Item(2,3).Value = 10
$totals.Cells.Item(2,4).Value = 4.36

# Renumber the pre-existing rows' index column (0..4 -> 1..5)
$totals.Cells.Item(3,1).Value = 1
$totals.Cells.Item(4,1).Value = 2
$totals.Cells.Item(5,1).Value = 3
$totals.Cells.Item(6,1).Value = 4
$totals.Cells.Item(7,1).Value = 5
